$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Canje_de_Puntos" (club) cases and their associated method codes.
$ws.Range("A17").Value = "Canje_de_Puntos_Canje_de_Credito_MIX"
$ws.Range("B17").Value = 1162816939
$ws.Range("B17").HorizontalAlignment = -4152

$ws.Range("A18").Value = "Canje_de_Puntos_Canje_de_Credito_PRE"
$ws.Range("B18").Value = 1162676705

$ws.Range("A19").Value = "Canje_de_Puntos_Canje_de_Pack_MIX"
$ws.Range("B19").Value = 1162816939
$ws.Range("B19").HorizontalAlignment = -4152

$ws.Range("A20").Value = "Canje_de_Puntos_Canje_de_Pack_PRE"
$ws.Range("B20").Value = 1162676705

$ws.Range("D21").Select() | Out-Null
